# Add Test Data for Hungary/Russia/Finland Market
# Each new sheet is modeled on the existing "Denmark" sheet (same layout,
# merged header cells, column widths and row styling), with the country
# specific cells (market name + NGC code) updated and, where the source
# data doesn't include the "MZXSDR240" repeater row, that row removed.

$wb = $excel.ActiveWorkbook
$denmark = $wb.Worksheets.Item("Denmark")

# ---------------------------------------------------------------------
# Russia - copied from Denmark, MZXSDR240 row removed
# ---------------------------------------------------------------------
$denmark.Copy($null, $denmark)
$russia = $wb.Worksheets.Item($denmark.Index + 1)
$russia.Name = "Russia"
$russia.Rows.Item(16).Delete()
$russia.Range("B2").Value = "Russia Market"
$russia.Range("B4").Value = "NGC-2929/T2900"
$russia.Rows.Item("3:5").RowHeight = 28.8

# ---------------------------------------------------------------------
# Finland - copied from Denmark, keeps the MZXSDR240 row
# ---------------------------------------------------------------------
$denmark.Copy($null, $russia)
$finland = $wb.Worksheets.Item($russia.Index + 1)
$finland.Name = "Finland"
$finland.Range("B2").Value = "Finland Market"
$finland.Range("B4").Value = "NGC-3130/T2943"
$finland.Rows.Item("3:5").RowHeight = 28.8

# ---------------------------------------------------------------------
# Hungary - copied from Denmark, MZXSDR240 row removed
# ---------------------------------------------------------------------
$denmark.Copy($null, $finland)
$hungary = $wb.Worksheets.Item($finland.Index + 1)
$hungary.Name = "Hungary"
$hungary.Rows.Item(16).Delete()
$hungary.Range("B2").Value = "Hungary Market"
$hungary.Range("B4").Value = "NGC-3104/T2992"
$hungary.Rows.Item("3:5").RowHeight = 28.8

# Hungary is the last sheet added, so it becomes the active / visible tab
$hungary.Activate()
